$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SkillBase")

# ---------------------------------------------------------------------------
# 1. Propagate correct cell formatting (style) onto the rows that are being
#    turned into real data rows. Rows 23-28 already carry the right
#    alternating banding style (they were blank template rows). Row 29
#    currently carries a stale "last row" style and rows 30-34 don't exist
#    yet, so copy the format of the last two fully-styled data rows
#    (21:22 - one "odd" banded row + one "even" banded row, which also
#    carries the special unbanded style on its ACT-seconds column) down
#    across the new rows.
# ---------------------------------------------------------------------------
$ws.Range("C21:AC22").Copy() | Out-Null
$ws.Range("C29:AC30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C31:AC32").PasteSpecial(-4122) | Out-Null
$ws.Range("C33:AC34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the new skill rows (105001-105004, 106001-106004, 107001-107004)
# ---------------------------------------------------------------------------
$rows = @(
  @{ r=23; id=105001; lvl=1; typ=1; atb=5000;  rage=10;   mp=0;  fcd=0; gcd=0;    lim=1; rng=2;  rt=1; tt=4; tl=0; tw=0; sc=1; eff=1010011 },
  @{ r=24; id=105002; lvl=1; typ=2; atb=10000; rage=20;   mp=0;  fcd=0; gcd=1000; lim=1; rng=10; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=25; id=105003; lvl=1; typ=3; atb=20000; rage=-100; mp=0;  fcd=0; gcd=2000; lim=1; rng=40; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=26; id=105004; lvl=1; typ=4; atb=2000;  rage=0;    mp=20; fcd=0; gcd=1500; lim=1; rng=20; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=27; id=106001; lvl=1; typ=1; atb=5000;  rage=10;   mp=0;  fcd=0; gcd=0;    lim=1; rng=12; rt=1; tt=4; tl=0; tw=0; sc=1; eff=1010011 },
  @{ r=28; id=106002; lvl=1; typ=2; atb=10000; rage=20;   mp=0;  fcd=0; gcd=1000; lim=1; rng=16; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=29; id=106003; lvl=1; typ=3; atb=20000; rage=-100; mp=0;  fcd=0; gcd=2000; lim=1; rng=40; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=30; id=106004; lvl=1; typ=4; atb=2000;  rage=0;    mp=20; fcd=0; gcd=1500; lim=1; rng=20; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=31; id=107001; lvl=1; typ=1; atb=5000;  rage=10;   mp=0;  fcd=0; gcd=0;    lim=1; rng=12; rt=1; tt=4; tl=0; tw=0; sc=1; eff=1010011 },
  @{ r=32; id=107002; lvl=1; typ=2; atb=10000; rage=20;   mp=0;  fcd=0; gcd=1000; lim=1; rng=16; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=33; id=107003; lvl=1; typ=3; atb=20000; rage=-100; mp=0;  fcd=0; gcd=2000; lim=1; rng=40; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 },
  @{ r=34; id=107004; lvl=1; typ=4; atb=2000;  rage=0;    mp=20; fcd=0; gcd=1500; lim=1; rng=20; rt=1; tt=4; tl=4; tw=0; sc=5; eff=1010011 }
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Range("C$r").Value = $row.id
  $ws.Range("D$r").Value = $row.lvl
  $ws.Range("M$r").Value = $row.typ
  $ws.Range("N$r").Value = $row.atb
  if ($r -lt 27) {
    $ws.Range("O$r").Value = 10000 / $row.atb
  }
  $ws.Range("P$r").Value = $row.rage
  $ws.Range("Q$r").Value = $row.mp
  $ws.Range("R$r").Value = $row.fcd
  $ws.Range("S$r").Value = $row.gcd
  $ws.Range("T$r").Value = $row.lim
  $ws.Range("U$r").Value = $row.rng
  $ws.Range("V$r").Value = $row.rt
  $ws.Range("W$r").Value = $row.tt
  $ws.Range("X$r").Value = $row.tl
  $ws.Range("Y$r").Value = $row.tw
  $ws.Range("Z$r").Value = $row.sc
  $ws.Range("AA$r").Value = $row.eff
}

# Rows 27-34 share one calculated formula (mirrors the earlier O8:O22 shared
# group) - assign across the whole block at once so the writer emits it as a
# single shared formula the way Excel does when a table's calculated column
# is extended.
$ws.Range("O27:O34").Formula = "=10000/N27"

# ---------------------------------------------------------------------------
# 3. Grow the "表5" table to cover the newly populated rows
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("C1:AC34"))

# ---------------------------------------------------------------------------
# 4. Restore the cursor/selection state
# ---------------------------------------------------------------------------
$ws.Range("I22").Select()
